$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 22 updates (@CrazyProgrammer_IT_IS)
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = "2026-02-19T02:37:01.845711+00:00"
$ws.Range("I22").Value = 1
$ws.Range("M22").Value = "[32702]"
